# Reordering exercises within each day on the "Workout Plan" sheet.
# No exercises were added/removed/changed in type - only their order
# (rows) within a given day was swapped. Columns A (Day) and B (Order)
# stay put; the exercise details in C:G are exchanged between the
# row-pairs below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workout Plan")

function Swap-Rows($ws, $r1, $r2) {
    $range1 = $ws.Range("C" + $r1 + ":G" + $r1)
    $range2 = $ws.Range("C" + $r2 + ":G" + $r2)
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# Day 1
Swap-Rows $ws 2 3
Swap-Rows $ws 4 5

# Day 2
Swap-Rows $ws 12 13
Swap-Rows $ws 15 16

# Day 3
Swap-Rows $ws 18 19
Swap-Rows $ws 21 22

# Day 4
Swap-Rows $ws 26 27
Swap-Rows $ws 28 31

# Update the last active selection/cell to match the edited workbook's
# view state.
$ws.Range("J35").Select()
